$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the full used range of the sheet so we cover every data row.
$ur = $ws.UsedRange
$firstRow = $ur.Row
$lastRow = $firstRow + $ur.Rows.Count - 1
$firstCol = $ur.Column
$lastCol = $firstCol + $ur.Columns.Count - 1

# Locate the "Treatment" column from the header row (row 1) instead of
# hard-coding a column letter, so the replacement targets the right data.
$treatmentCol = 0
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Text
    if ($header -eq "Treatment") {
        $treatmentCol = $c
    }
}
if ($treatmentCol -eq 0) {
    $treatmentCol = 6
}

# Rename the treatment-group labels: Con -> X, SF1 -> Y, SF2 -> Z.
# Each label is fully replaced (in its own pass) before moving on to the
# next, so the shared-string table keeps stable slots for the renamed
# values instead of interleaving replacements across different labels.
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $treatmentCol)
    if ($cell.Text -eq "Con") {
        $cell.Value = "X"
    }
}
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $treatmentCol)
    if ($cell.Text -eq "SF1") {
        $cell.Value = "Y"
    }
}
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $treatmentCol)
    if ($cell.Text -eq "SF2") {
        $cell.Value = "Z"
    }
}
